# Generate Report for Handback
# Regenerates the handback-status report timestamps/priority that were
# produced for the 343aca44-... and cf740957-... entries (these two rows
# share identical values throughout the workbook, so both move together).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview.Range("G2").Value = "2016-08-15 18:15:33"
$wsOverview.Range("G4").Value = "2016-08-15 18:15:33"

# --- zh-cn sheet ---
# Priority (column E)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime (column H)
$wsZhCn.Range("H2").Value = "2016-08-15 18:15:28"
$wsZhCn.Range("H4").Value = "2016-08-15 18:15:28"
# Correspond Handback DateTime (column K)
$wsZhCn.Range("K2").Value = "2016-08-15 18:15:44"
$wsZhCn.Range("K4").Value = "2016-08-15 18:15:44"

# --- de-de sheet ---
# Priority (column E)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime (column H)
$wsDeDe.Range("H2").Value = "2016-08-15 18:15:51"
$wsDeDe.Range("H4").Value = "2016-08-15 18:15:51"
